# Update saved 2D analytic-fit parameters (sheet "Analitico"), add a new
# "RMSE" column to both sheets, and refresh the plots/selection/active-tab
# state to match the latest run ("update plots and saved 2d analytic
# parameters").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Analitico
$ws2 = $wb.Worksheets.Item(2)   # Numerico

function Set-Row {
    param($ws, [int]$row, [string[]]$vals)
    $cols = @("B", "C", "D", "E", "F", "G", "H", "I")
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($vals[$i] -ne $null) {
            $ws.Range("$($cols[$i])$row").Value = [double]$vals[$i]
        }
    }
}

# --- New "RMSE" header (column I) on both sheets ------------------------
$ws1.Range("I1").Value = "RMSE"
$ws2.Range("I1").Value = "RMSE"

# --- Analitico: refreshed fit results (columns B:I, rows 2-12) ----------
# (Full-precision literals so the round-tripped IEEE-754 double exactly
# matches the value Excel originally stored.)
Set-Row $ws1 2  @("0.76100000000000001", "2.48E-6",                "-6.81E-6",                "3.5400000000000001E-2",  "45.2",  "1.67",                "2",    "0.45")
Set-Row $ws1 3  @("0.52400000000000002", "2.2E-12",                "-3.3700000000000001E-9",  "7.3999999999999996E-2",  "227",   "1.25",                "2",    "7.26")
Set-Row $ws1 4  @("0.46300000000000002", "1.4700000000000002E-11", "-2.1200000000000001E-9",  "2.3400000000000001E-2",  "344",   "1.44",                "2",    "1.22")
Set-Row $ws1 5  @("0.52",                "3.0400000000000002E-14", "-1.38E-9",                "9.7099999999999999E-3",  "4780",  "1.1299999999999999",  "2",    "1.98")
Set-Row $ws1 6  @("1.03",                "4.0800000000000002E-5",  "-1.7000000000000001E-4",  "1.5",                     "686",   "1.58",                "2",    "1.61")
Set-Row $ws1 7  @("8.2200000000000006",  "2.2800000000000001E-4",  "-3.6499999999999998E-4",  "0.193",                   "127",   "1.9",                 "2",    "2.5499999999999998")
Set-Row $ws1 8  @("0.504",                "3.5700000000000001E-10","-1.26E-8",                 "0.109",                   "775",   "1.63",                "2",    "2.76")
Set-Row $ws1 9  @("8.2200000000000006",  "4.2899999999999999E-5",  "1.8800000000000001E-2",   "1.6500000000000001E-2",  "0.188", "4.7199999999999999E-2","1.8", "4.68")
Set-Row $ws1 10 @("0.47299999999999998", "1.48E-12",               "-4.6000000000000001E-10", "6.6600000000000006E-2",  "66800", "1.28",                "2",    "4.6100000000000003")
Set-Row $ws1 11 @("0.45400000000000001", "9.1900000000000002E-13", "-3.4200000000000002E-9",  "0.11799999999999999",    "240",   "1.23",                "2",    "3.3")
Set-Row $ws1 12 @("0.46700000000000003", "-2.4600000000000001E-43","7.3199999999999994E-38",  "0.89900000000000002",    "6290",  "1.08",                "1.24", "1.22")

# --- Numerico: only the new RMSE column is populated (rows 2-8) ---------
$ws2.Range("I2").Value = [double]"0.70799999999999996"
$ws2.Range("I3").Value = [double]"2.78"
$ws2.Range("I4").Value = [double]"1.1499999999999999"
$ws2.Range("I5").Value = [double]"0.21"
$ws2.Range("I6").Value = [double]"0.32200000000000001"
$ws2.Range("I7").Value = [double]"1.66"
$ws2.Range("I8").Value = [double]"0.76800000000000002"

# --- Selection / active-sheet bookkeeping --------------------------------
# Analitico is no longer the selected tab; its lingering selection moves to L12.
$ws1.Range("L12").Select()

# Numerico becomes the active / selected tab, with its selection at K7.
$ws2.Activate()
$ws2.Range("K7").Select()
